$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.988.14'
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").Value = '2.373.43'
$ws.Range("E3").Value = '  +0.78%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''0.691'
$ws.Range("E5").Value = '  +6.30%  '

$ws.Range("D6").Value = '''242.56'
$ws.Range("E6").Value = '  +3.28%  '

$ws.Range("D7").Value = '''76.68'
$ws.Range("E7").Value = '  +7.17%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.638'
$ws.Range("E9").Value = '  +28.88%  '

$ws.Range("E10").Value = '  +5.27%  '

$ws.Range("D11").Value = '''57.56'
$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").Value = '''33.18'
$ws.Range("E12").Value = '  +21.32%  '

$ws.Range("D13").Value = '''7.62'
$ws.Range("E13").Value = '  +19.97%  '

$ws.Range("E14").Value = '  +1.86%  '

$ws.Range("D15").Value = '2.727.54'
$ws.Range("E15").Value = '  +0.73%  '

$ws.Range("D16").Value = '''16.98'
$ws.Range("E16").Value = '  +4.77%  '

$ws.Range("E17").Value = '  +7.13%  '

$ws.Range("D18").Value = '2.370.99'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '44.079.62'
$ws.Range("E19").Value = '  +1.92%  '

$ws.Range("E20").Value = '  +2.77%  '

$ws.Range("D21").Value = '''6.71'
$ws.Range("E21").Value = '  +5.68%  '

$ws.Range("D22").Value = '''78.08'
$ws.Range("E22").Value = '  +4.00%  '

$ws.Range("D23").Value = '''258.80'
$ws.Range("E23").Value = '  +3.33%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '''2.55'
$ws.Range("E25").Value = '  +3.53%  '

$ws.Range("E26").Value = '  +11.28%  '

$ws.Range("E27").Value = '  -2.28%  '

$ws.Range("D28").Value = '''1.77'
$ws.Range("E28").Value = '  +15.35%  '

$ws.Range("E29").Value = '  +3.77%  '

$ws.Range("D30").Value = '''23.28'
$ws.Range("E30").Value = '  +3.92%  '

$ws.Range("D31").Value = '''176.23'
$ws.Range("E31").Value = '  +2.21%  '

$ws.Range("E32").Value = '  -0.41%  '

$ws.Range("E33").Value = '  +5.59%  '

$ws.Range("D34").Value = '''5.37'
$ws.Range("E34").Value = '  +7.24%  '

$ws.Range("E35").Value = '  +8.87%  '

$ws.Range("E36").Value = '  +6.53%  '

$ws.Range("D37").Value = '''3.84'
$ws.Range("E37").Value = '  +2.67%  '

$ws.Range("D38").Value = '''2.47'
$ws.Range("E38").Value = '  +1.51%  '

$ws.Range("D39").Value = '''6.53'
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("D40").Value = '''0.0277'
$ws.Range("E40").Value = '  +8.91%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '''0.204'
$ws.Range("E41").Value = '  +19.95%  '

$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").Value = '''19.11'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("D43").Value = '''9.02'
$ws.Range("E43").Value = '  +1.12%  '

$ws.Range("E44").Value = '  +0.00%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = '''0.101'
$ws.Range("E45").Value = '  +5.48%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''1.21'
$ws.Range("E46").Value = '  +4.34%  '

$ws.Range("E47").Value = '  +5.02%  '

$ws.Range("E48").Value = '  +14.13%  '

$ws.Range("D49").Value = '''103.22'
$ws.Range("E49").Value = '  +4.37%  '

$ws.Range("E50").Value = '  +0.23%  '

$ws.Range("D51").Value = '''54.85'
$ws.Range("E51").Value = '  +8.47%  '
